$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update R-column delay values (numeric)
$ws.Range("R2").Value = -0.3623372395833334
$ws.Range("R3").Value = -0.4377821180555556
$ws.Range("R4").Value = -0.4741319444444445
$ws.Range("R6").Value = 0
$ws.Range("R7").Value = -0.5925130208333333
$ws.Range("R8").Value = -1.416059027777778
$ws.Range("R9").Value = -1.443691314548611
$ws.Range("R10").Value = -1.406147540983796
$ws.Range("R11").Value = 0
$ws.Range("R12").Value = -1.498667800451389
$ws.Range("R13").Value = -0.4885587431712963
$ws.Range("R14").Value = -1.514301215277778
$ws.Range("R15").Value = -0.5466647104861111
$ws.Range("R16").Value = -2.611839708564815
$ws.Range("R17").Value = -2.471889671365741
$ws.Range("R18").Value = -0.5261067708333333
$ws.Range("R19").Value = -2.5365234375
$ws.Range("R20").Value = 0
$ws.Range("R22").Value = -0.4133138020833333
$ws.Range("R23").Value = -1.526104797974537
$ws.Range("R24").Value = -1.442361111111111
$ws.Range("R25").Value = -1.466666666666667

# Update N-column "veicolo" cells to annotated text values
$ws.Range("N5").Value = "39666 (non in estrazione)"
$ws.Range("N6").Value = "39742 (esterno)"
$ws.Range("N11").Value = "39666 (esterno)"
$ws.Range("N20").Value = "39762 (esterno)"
$ws.Range("N21").Value = "39723 (non in estrazione)"
$ws.Range("N26").Value = "39750 (non in estrazione)"
$ws.Range("N27").Value = "39764 (non in estrazione)"
